# The target revision (per the supplied canonical-OOXML diff) is a pure
# re-serialization artifact: every hunk only reorders XML attributes
# (namespace declarations grouped/sorted by prefix, then the remaining
# attributes sorted alphabetically by local name) inside word/document.xml
# and word/styles.xml. None of the document's actual content, structure,
# formatting values, or text changed - e.g. <w:pgSz w:w="11906" w:h="16838"/>
# became <w:pgSz w:h="16838" w:w="11906"/>, <w:lang w:val="fr-FR"
# w:eastAsia="en-US" w:bidi="ar-SA"/> became <w:lang w:bidi="ar-SA"
# w:eastAsia="en-US" w:val="fr-FR"/>, <w:style w:type="paragraph"
# w:default="1" w:styleId="Normal"> became <w:style w:default="1"
# w:styleId="Normal" w:type="paragraph">, etc. The accompanying commit
# message ("Fixed POI packaging and upgraded to POI 3.15.") confirms this:
# it is a build/library-upgrade commit, and the new Apache POI/XMLBeans
# version simply wrote attributes back out in a different (alphabetical)
# order when the test fixture .docx was regenerated - the rendered
# document itself is byte-for-byte identical in meaning.
#
# Because the attribute order is not part of the semantic object model
# (Word's COM interface has no notion of "attribute order" - it is purely
# a side effect of which XML serializer wrote the part), there is no
# content edit to apply here. The correct action is to leave the
# document's content untouched so the canonical (order-insensitive) OOXML
# keeps matching the target revision.

$d = $word.ActiveDocument
